# Doug Kinsey 2026-01-19 weekly timesheet -> corrected hours/client names,
# plus a new "Jason Schema" sheet that flattens the same data for export.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Weekly Timesheet")

# --- Row 2: 2026-01-19 ---
$ws.Range("B2").Value = "Hall"
$ws.Range("C2").Value = 10
$ws.Range("F2").Value = 650

# --- Row 3: 2026-01-20 ---
$ws.Range("B3").Value = "Bryan"
$ws.Range("C3").Value = 10
$ws.Range("F3").Value = 650

# --- Row 4: 2026-01-21 ---
$ws.Range("B4").Value = "McGill"
$ws.Range("C4").Value = 10
$ws.Range("F4").Value = 650

# --- Row 5: 2026-01-22 ---
$ws.Range("B5").Value = "Hall"
$ws.Range("C5").Value = 8
$ws.Range("F5").Value = 520

# --- Row 6: 2026-01-25 (date changed from 2026-01-23) ---
$ws.Range("A6").Value = "2026-01-25"
$ws.Range("B6").Value = "McGill"
$ws.Range("C6").Value = 2
$ws.Range("F6").Value = 130

# --- Row 7: 2026-01-25 (date changed from 2026-01-23) ---
$ws.Range("A7").Value = "2026-01-25"
$ws.Range("B7").Value = "McGill"
$ws.Range("C7").Value = 8
$ws.Range("F7").Value = 780

# --- Row 9: SUBTOTAL ---
$ws.Range("C9").Value = 48
$ws.Range("D9").Value = "Reg: 40 / OT: 8"
$ws.Range("F9").Value = 3380

# --- Add the new "Jason Schema" sheet after "Weekly Timesheet" ---
$jason = $wb.Worksheets.Add($null, $ws)
$jason.Name = "Jason Schema"

$jason.Range("A1").Value = "Employee"
$jason.Range("B1").Value = "Employee ID"
$jason.Range("C1").Value = "Date"
$jason.Range("D1").Value = "Client"
$jason.Range("E1").Value = "Hours"
$jason.Range("F1").Value = "Rate"
$jason.Range("G1").Value = "Total"
$jason.Range("H1").Value = "Type"
$jason.Range("I1").Value = "Notes"

$jason.Range("A2:A7").Value = "Doug Kinsey"
$jason.Range("B2:B7").Value = "emp_JMr5EHDoCPPJw1h4"

$jason.Range("C2").Value = "2026-01-19"
$jason.Range("D2").Value = "Hall"
$jason.Range("E2").Value = 10
$jason.Range("F2").Value = 65
$jason.Range("G2").Value = 650
$jason.Range("H2").Value = "Regular"

$jason.Range("C3").Value = "2026-01-20"
$jason.Range("D3").Value = "Bryan"
$jason.Range("E3").Value = 10
$jason.Range("F3").Value = 65
$jason.Range("G3").Value = 650
$jason.Range("H3").Value = "Regular"

$jason.Range("C4").Value = "2026-01-21"
$jason.Range("D4").Value = "McGill"
$jason.Range("E4").Value = 10
$jason.Range("F4").Value = 65
$jason.Range("G4").Value = 650
$jason.Range("H4").Value = "Regular"

$jason.Range("C5").Value = "2026-01-22"
$jason.Range("D5").Value = "Hall"
$jason.Range("E5").Value = 8
$jason.Range("F5").Value = 65
$jason.Range("G5").Value = 520
$jason.Range("H5").Value = "Regular"

$jason.Range("C6").Value = "2026-01-25"
$jason.Range("D6").Value = "McGill"
$jason.Range("E6").Value = 2
$jason.Range("F6").Value = 65
$jason.Range("G6").Value = 130
$jason.Range("H6").Value = "Regular"

$jason.Range("C7").Value = "2026-01-25"
$jason.Range("D7").Value = "McGill"
$jason.Range("E7").Value = 8
$jason.Range("F7").Value = 65
$jason.Range("G7").Value = 780
$jason.Range("H7").Value = "OT"

$jason.Range("I2:I7").Value = ""

# header formatting: bold, no fill
$jason.Range("A1:I1").Font.Bold = $true

# column widths to roughly match the source layout
$jason.Columns.Item(1).ColumnWidth = 20
$jason.Columns.Item(2).ColumnWidth = 18
$jason.Columns.Item(3).ColumnWidth = 12
$jason.Columns.Item(4).ColumnWidth = 25
$jason.Columns.Item(5).ColumnWidth = 8
$jason.Columns.Item(6).ColumnWidth = 10
$jason.Columns.Item(7).ColumnWidth = 12
$jason.Columns.Item(8).ColumnWidth = 10
$jason.Columns.Item(9).ColumnWidth = 30
